$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Program")
$ws.Activate()

# Update the test data values in place (mirrors the "clean run" refresh of
# the randomized test data used by the LMS automation suite).
$ws.Range("A5").Value = "Team4kl8Team4Team4"
$ws.Range("B5").Value = "javajavkli09Team4"
$ws.Range("A6").Value = "JavakjhgflTeam4"
$ws.Range("B7").Value = "XkjmnhjkmnhjTeam4"

# Move the active selection from A5 to B7.
$ws.Range("B7").Select()
